$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AW2").Value = 162.230868
$ws.Range("AQ3").Value = 25.02022
$ws.Range("AW4").Value = 5.346863
$ws.Range("AW5").Value = 117.176528
$ws.Range("AW6").Value = 119.963079
$ws.Range("AK7").Value = 20.068194
$ws.Range("AW8").Value = 155.000579
$ws.Range("AW9").Value = 104.041215
$ws.Range("AW10").Value = 88.12450200000001
$ws.Range("AW11").Value = 4.021146
$ws.Range("AW12").Value = 89.026875
$ws.Range("AW13").Value = 167.084572
$ws.Range("AW14").Value = 152.899722
$ws.Range("AW15").Value = 122.028252
$ws.Range("AW16").Value = 124.325556
$ws.Range("AN17").Value = 27.198588
$ws.Range("AW18").Value = 3.97537
$ws.Range("AW19").Value = 91.912384
$ws.Range("AW20").Value = 84.00114600000001
$ws.Range("AK21").Value = 20.068206
$ws.Range("AK22").Value = 20.068252
$ws.Range("AW23").Value = 5.347292
$ws.Range("AW24").Value = 117.175995
$ws.Range("AW25").Value = 4.120208
$ws.Range("AK26").Value = 20.068356
$ws.Range("AW27").Value = 152.905231
$ws.Range("AW28").Value = 83.082199
$ws.Range("AW29").Value = 73.988715
$ws.Range("AW30").Value = 13.222245
$ws.Range("AW31").Value = 154.995486
$ws.Range("AW32").Value = 91.22420099999999
$ws.Range("AW33").Value = 12.003646
$ws.Range("AW34").Value = 0.121991
$ws.Range("AW35").Value = 119.27
$ws.Range("AW36").Value = 127.055741
$ws.Range("AK37").Value = 20.068229
$ws.Range("AW38").Value = 14.287639
$ws.Range("AK39").Value = 20.068438
$ws.Range("AK40").Value = 20.068148
$ws.Range("AK41").Value = 20.068125
$ws.Range("AQ42").Value = 25.020093
$ws.Range("AW43").Value = 76.307002
$ws.Range("AW44").Value = 91.22358800000001
$ws.Range("AW45").Value = 42.340359
$ws.Range("AW46").Value = 42.339144
$ws.Range("AK47").Value = 42.117743
$ws.Range("AW48").Value = 119.270984
$ws.Range("AW49").Value = 119.964491
$ws.Range("AW50").Value = 27.107106
$ws.Range("AW51").Value = 133.179977
$ws.Range("AW52").Value = 110.937512
$ws.Range("AW53").Value = 158.042164
$ws.Range("AW54").Value = 148.090984
$ws.Range("AK55").Value = 42.917986
$ws.Range("AW56").Value = 119.963657
$ws.Range("AQ57").Value = 25.019606
$ws.Range("AW58").Value = 14.291748
$ws.Range("AW59").Value = 123.21934
$ws.Range("AW60").Value = 84.02048600000001
$ws.Range("AW61").Value = 55.965231
$ws.Range("AK62").Value = 15.093958
$ws.Range("AW63").Value = 4.020359
$ws.Range("AW64").Value = 159.335301
$ws.Range("AW65").Value = 162.228194
$ws.Range("AW66").Value = 119.979433
$ws.Range("AW67").Value = 77.080961
$ws.Range("AW68").Value = 20.113519
$ws.Range("AK69").Value = 20.068171
$ws.Range("AK70").Value = 20.068414
$ws.Range("AW71").Value = 13.117558
$ws.Range("AK72").Value = 42.117824
$ws.Range("AW73").Value = 134.050868
$ws.Range("AW74").Value = 71.077083
$ws.Range("AW75").Value = 6.004074
$ws.Range("AW76").Value = 75.285729
$ws.Range("AN77").Value = 13.233727
$ws.Range("AW78").Value = 148.094375
$ws.Range("AW79").Value = 27.255498
$ws.Range("AW80").Value = 42.340139
$ws.Range("AW81").Value = 21.950255
$ws.Range("AW82").Value = 12.002303
$ws.Range("AW83").Value = 127.119699
$ws.Range("AW84").Value = 157.994444
$ws.Range("AW85").Value = 105.964109
$ws.Range("AW86").Value = 75.28401599999999
$ws.Range("AW87").Value = 55.958079
$ws.Range("AQ88").Value = 40.049884
$ws.Range("AW89").Value = 111.96662
$ws.Range("AW90").Value = 152.89934
$ws.Range("AW91").Value = 62.339769
$ws.Range("AK92").Value = 15.091516
$ws.Range("AW93").Value = 119.962882
$ws.Range("AK94").Value = 20.068079
$ws.Range("AW95").Value = 159.335694
$ws.Range("AW96").Value = 70.95035900000001
$ws.Range("AQ97").Value = 25.020324
$ws.Range("AW98").Value = 89.019537
$ws.Range("AW99").Value = 81.095394
$ws.Range("AW100").Value = 155.085301
$ws.Range("AW101").Value = 111.082292
$ws.Range("AW102").Value = 50.341042
$ws.Range("AK103").Value = 15.093981
$ws.Range("AQ104").Value = 8.246123000000001
$ws.Range("AW105").Value = 71.006308
$ws.Range("AW106").Value = 109.519167
$ws.Range("AW107").Value = 118.162685
$ws.Range("AK108").Value = 20.068113
$ws.Range("AW109").Value = 20.095
$ws.Range("AN110").Value = 13.233681
$ws.Range("AK111").Value = 20.068299
$ws.Range("AW112").Value = 162.231157
$ws.Range("AW113").Value = 157.992581
$ws.Range("AW114").Value = 148.042037
$ws.Range("AW115").Value = 50.339444
$ws.Range("AW116").Value = 4.119722
$ws.Range("AW117").Value = 21.940891
$ws.Range("AW118").Value = 88.97274299999999
$ws.Range("AW119").Value = 136.96125
$ws.Range("AW120").Value = 111.102303
$ws.Range("AW121").Value = 75.281863
$ws.Range("AK122").Value = 20.068241
$ws.Range("AK123").Value = 42.91809
$ws.Range("AW124").Value = 152.899954
$ws.Range("AK125").Value = 20.068148
$ws.Range("AN126").Value = 13.23375
$ws.Range("AW127").Value = 147.994549
$ws.Range("AK128").Value = 42.117859
$ws.Range("AK129").Value = 20.068125
$ws.Range("AW130").Value = 62.202118
$ws.Range("AW131").Value = 56.35353
$ws.Range("AW132").Value = 166.044572
$ws.Range("AW133").Value = 91.88493099999999
$ws.Range("AW134").Value = 50.3389
$ws.Range("AN135").Value = 27.123461
$ws.Range("AK136").Value = 20.068218
$ws.Range("AW137").Value = 148.045243
$ws.Range("AW138").Value = 152.901609
$ws.Range("AW139").Value = 89.082234
$ws.Range("AK140").Value = 20.06816
$ws.Range("AW141").Value = 123.218912
$ws.Range("AW142").Value = 119.96338
$ws.Range("AW143").Value = 83.083229
$ws.Range("AK144").Value = 20.068403
$ws.Range("AW145").Value = 55.015463
$ws.Range("AW146").Value = 55.215405
$ws.Range("AQ147").Value = 40.049769
$ws.Range("AW148").Value = 76.9136
$ws.Range("AW149").Value = 66.978831
$ws.Range("AQ150").Value = 32.029884
$ws.Range("AW151").Value = 97.04900499999999
$ws.Range("AW152").Value = 50.338206
$ws.Range("AK153").Value = 20.068264
$ws.Range("AW154").Value = 119.270787
$ws.Range("AW155").Value = 75.287477
$ws.Range("AK156").Value = 15.091516
$ws.Range("AN157").Value = 8.100984
$ws.Range("AW158").Value = 109.2936
$ws.Range("AW159").Value = 25.040046
$ws.Range("AQ160").Value = 3.938125
$ws.Range("AW161").Value = 157.995382
$ws.Range("AW162").Value = 134.345683
$ws.Range("AW163").Value = 73.988507
$ws.Range("AK164").Value = 20.06809
$ws.Range("AW165").Value = 109.298542
$ws.Range("AW166").Value = 162.231551
$ws.Range("AW167").Value = 83.999618
$ws.Range("AW168").Value = 154.994931
$ws.Range("AW169").Value = 146.101771
$ws.Range("AW170").Value = 84.02037
$ws.Range("AW171").Value = 13.117396
$ws.Range("AK172").Value = 42.117731
$ws.Range("AW173").Value = 154.994572
$ws.Range("AW174").Value = 111.098113
$ws.Range("AW175").Value = 56.352813
$ws.Range("AW176").Value = 88.91252299999999
$ws.Range("AW177").Value = 0.065868
$ws.Range("AW178").Value = 162.23912
$ws.Range("AW179").Value = 97.06119200000001
$ws.Range("AW180").Value = 75.284988
$ws.Range("AW181").Value = 56.356354
$ws.Range("AK182").Value = 20.068241
$ws.Range("AW183").Value = 111.101933
$ws.Range("AK184").Value = 42.117859
$ws.Range("AW185").Value = 169.306412
$ws.Range("AW186").Value = 155.085579
$ws.Range("AW187").Value = 75.29083300000001
$ws.Range("AW188").Value = 49.078229
$ws.Range("AW189").Value = 142.322118
$ws.Range("AW190").Value = 152.900521
$ws.Range("AW191").Value = 81.004074
$ws.Range("AW192").Value = 136.95809
$ws.Range("AW193").Value = 84.92751199999999
$ws.Range("AK194").Value = 15.09397
$ws.Range("AW195").Value = 157.968715
$ws.Range("AW196").Value = 155.085035
$ws.Range("AW197").Value = 168.130231
$ws.Range("AW198").Value = 154.007257
$ws.Range("AW199").Value = 165.989838
$ws.Range("AW200").Value = 158.043414
$ws.Range("AW201").Value = 119.965104
$ws.Range("AW202").Value = 25.308113
$ws.Range("AW203").Value = 119.983079
$ws.Range("AK204").Value = 42.117743
$ws.Range("AW205").Value = 167.085556
$ws.Range("AW206").Value = 157.984479
$ws.Range("AW207").Value = 167.081424
$ws.Range("AW208").Value = 68.035729
$ws.Range("AK209").Value = 20.068194
$ws.Range("AK210").Value = 20.068183
$ws.Range("AQ211").Value = 25.019826
$ws.Range("AW212").Value = 0.09002300000000001
$ws.Range("AW213").Value = 75.29018499999999
$ws.Range("AW214").Value = 127.05566
$ws.Range("AW215").Value = 144.296968
$ws.Range("AW216").Value = 84.066574
$ws.Range("AW217").Value = 76.306771
$ws.Range("AW218").Value = 56.354664
$ws.Range("AW219").Value = 7.28353
$ws.Range("AW220").Value = 119.271319
$ws.Range("AW221").Value = 75.14978000000001
$ws.Range("AK222").Value = 42.117847
$ws.Range("AQ223").Value = 32.022523
$ws.Range("AW224").Value = 109.297315
$ws.Range("AW225").Value = 154.969815
$ws.Range("AW226").Value = 119.963889
$ws.Range("AW227").Value = 83.997963
$ws.Range("AW228").Value = 27.252546
$ws.Range("AW229").Value = 56.057431
$ws.Range("AW230").Value = 50.337326
$ws.Range("AK231").Value = 42.918102
$ws.Range("AK232").Value = 20.068137
$ws.Range("AW233").Value = 119.964294
$ws.Range("AW234").Value = 91.22482599999999
$ws.Range("AK235").Value = 20.068218
$ws.Range("AW236").Value = 140.117488
$ws.Range("AK237").Value = 20.061215
$ws.Range("AK238").Value = 20.06809
$ws.Range("AW239").Value = 119.231898
$ws.Range("AW240").Value = 50.344919
$ws.Range("AW241").Value = 50.347072
$ws.Range("AW242").Value = 24.929468
$ws.Range("AK243").Value = 20.068171
$ws.Range("AW244").Value = 117.176157
$ws.Range("AW245").Value = 41.335637
$ws.Range("AW246").Value = 27.107859
$ws.Range("AW247").Value = 148.094595
$ws.Range("AW248").Value = 158.042674
$ws.Range("AW249").Value = 136.958438
$ws.Range("AW250").Value = 27.108669
$ws.Range("AN251").Value = 13.233704
$ws.Range("AW252").Value = 13.117292
$ws.Range("AW253").Value = 136.958981
$ws.Range("AW254").Value = 56.349063
$ws.Range("AW255").Value = 8.095903
$ws.Range("AW256").Value = 146.987789
$ws.Range("AW257").Value = 62.20331
$ws.Range("AW258").Value = 14.289919
$ws.Range("AW259").Value = 133.17978
$ws.Range("AW260").Value = 13.222975
$ws.Range("AQ261").Value = 25.019734
$ws.Range("AW262").Value = 4.018981
$ws.Range("AW263").Value = 42.339063
$ws.Range("AQ264").Value = 7.174838
$ws.Range("AW265").Value = 134.346053
$ws.Range("AW266").Value = 146.10272
$ws.Range("AW267").Value = 89.026447
$ws.Range("AW268").Value = 109.29838
$ws.Range("AW269").Value = 30.353866
$ws.Range("AW270").Value = 133.308449
$ws.Range("AK271").Value = 42.117755
$ws.Range("AW272").Value = 8.089560000000001
$ws.Range("AW273").Value = 147.30765
$ws.Range("AW274").Value = 91.915544
$ws.Range("AW275").Value = 152.898553
$ws.Range("AW276").Value = 75.2914
$ws.Range("AW277").Value = 119.966227
$ws.Range("AQ278").Value = 4.023125
$ws.Range("AW279").Value = 144.297431
$ws.Range("AW280").Value = 75.284271
$ws.Range("AW281").Value = 50.095741
$ws.Range("AK282").Value = 20.068345
$ws.Range("AW283").Value = 159.337153
$ws.Range("AW284").Value = 85.057211
$ws.Range("AW285").Value = 49.009363
$ws.Range("AK286").Value = 20.068414
$ws.Range("AK287").Value = 20.068102
$ws.Range("AW288").Value = 4.119572
$ws.Range("AW289").Value = 49.009699
$ws.Range("AW290").Value = 49.007928
$ws.Range("AK291").Value = 42.918102
$ws.Range("AW292").Value = 152.898935
$ws.Range("AW293").Value = 105.960625
$ws.Range("AW294").Value = 25.040995
$ws.Range("AW295").Value = 105.963414
$ws.Range("AW296").Value = 119.268657
$ws.Range("AK297").Value = 15.091505
$ws.Range("AW298").Value = 162.232002
$ws.Range("AW299").Value = 56.348206
$ws.Range("AW300").Value = 20.058287
$ws.Range("AW301").Value = 8.093541999999999
$ws.Range("AW302").Value = 141.006053
$ws.Range("AK303").Value = 20.068322
$ws.Range("AK304").Value = 20.068275
$ws.Range("AW305").Value = 96.06506899999999
$ws.Range("AW306").Value = 119.965556
$ws.Range("AW307").Value = 88.972442
$ws.Range("AW308").Value = 34.98316
$ws.Range("AW309").Value = 14.286134
$ws.Range("AW310").Value = 4.01941
$ws.Range("AW311").Value = 98.99003500000001
$ws.Range("AK312").Value = 20.068345
$ws.Range("AW313").Value = 25.040509
$ws.Range("AW314").Value = 91.22147
$ws.Range("AW315").Value = 167.090382
$ws.Range("AW316").Value = 77.28935199999999
$ws.Range("AW317").Value = 56.355567
$ws.Range("AW318").Value = 6.003588
$ws.Range("AW319").Value = 124.328727
$ws.Range("AK320").Value = 20.068275
$ws.Range("AK321").Value = 20.068322
$ws.Range("AQ322").Value = 7.18338
$ws.Range("AW323").Value = 81.049109
$ws.Range("AW324").Value = 59.987627
$ws.Range("AW325").Value = 25.04162
$ws.Range("AW326").Value = 6.004815
$ws.Range("AW327").Value = 133.308669
$ws.Range("AW328").Value = 62.202859
$ws.Range("AW329").Value = 50.336563
$ws.Range("AW330").Value = 56.409074
$ws.Range("AK331").Value = 20.068391
$ws.Range("AQ332").Value = 8.246226999999999
$ws.Range("AW333").Value = 157.993287
$ws.Range("AW334").Value = 159.335463
$ws.Range("AK335").Value = 28.019722
$ws.Range("AK336").Value = 20.068287
$ws.Range("AW337").Value = 165.984132
$ws.Range("AW338").Value = 56.356782
$ws.Range("AW339").Value = 136.959502
$ws.Range("AW340").Value = 27.254051
$ws.Range("AW341").Value = 25.289086
$ws.Range("AW342").Value = 152.898287
$ws.Range("AW343").Value = 124.326042
$ws.Range("AW344").Value = 42.926238
$ws.Range("AW345").Value = 127.1186
$ws.Range("AW346").Value = 75.285625
$ws.Range("AW347").Value = 162.124039
$ws.Range("AW348").Value = 158.029896
$ws.Range("AW349").Value = 6.010208
$ws.Range("AW350").Value = 159.336262
$ws.Range("AW351").Value = 136.960023
$ws.Range("AW352").Value = 167.081782
$ws.Range("AW353").Value = 15.01816
$ws.Range("AW354").Value = 152.896817
$ws.Range("AW355").Value = 119.962697
$ws.Range("AW356").Value = 111.095324
$ws.Range("AW357").Value = 27.253368
$ws.Range("AK358").Value = 20.06831
$ws.Range("AK359").Value = 42.117836
$ws.Range("AK360").Value = 20.068426
$ws.Range("AW361").Value = 73.988854
$ws.Range("AW362").Value = 73.98839099999999
$ws.Range("AQ363").Value = 25.019977
$ws.Range("AW364").Value = 6.003796
